$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1302
$ws.Range("I19").Value = 1999.5
$ws.Range("J19").Value = 1102.7142
$ws.Range("K19").Value = 1999.5
$ws.Range("L19").Value = 1102.7142
$ws.Range("M19").Value = -1824.5
$ws.Range("N19").Value = -1452.7142
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = ""
$ws.Range("H58").Value = 943.5714
$ws.Range("I58").Value = 176.25
$ws.Range("K58").Value = 528.75
$ws.Range("M58").Value = -378.75
$ws.Range("H62").Value = 5831.967
$ws.Range("I62").Value = 5227.885
$ws.Range("J62").Value = 9758.5
$ws.Range("K62").Value = 5227.885
$ws.Range("L62").Value = 9758.5
$ws.Range("M62").Value = -4603.885
$ws.Range("N62").Value = -11006.5
$ws.Range("H65").Value = 5831.967
$ws.Range("I65").Value = 5227.885
$ws.Range("J65").Value = 9758.5
$ws.Range("K65").Value = 26139.425
$ws.Range("L65").Value = 48792.5
$ws.Range("M65").Value = -23019.425
$ws.Range("N65").Value = -55032.5
$ws.Range("H74").Value = 6177.4165
$ws.Range("I74").Value = 8089.857
$ws.Range("K74").Value = 8089.857
$ws.Range("M74").Value = -7153.857
$ws.Range("H77").Value = 6177.4165
$ws.Range("I77").Value = 8089.857
$ws.Range("K77").Value = 40449.285
$ws.Range("M77").Value = -35769.285
$ws.Range("H87").Value = 94155.336
$ws.Range("J87").Value = 79989
$ws.Range("L87").Value = 79989
$ws.Range("N87").Value = -82485
$ws.Range("H90").Value = 94155.336
$ws.Range("J90").Value = 79989
$ws.Range("L90").Value = 239967
$ws.Range("N90").Value = -252447
$ws.Range("H98").Value = 1435.8948
$ws.Range("I98").Value = 1413.25
$ws.Range("K98").Value = 1413.25
$ws.Range("M98").Value = 84.75
$ws.Range("H116").Value = 39001.25
$ws.Range("I116").Value = 39001.25
$ws.Range("K116").Value = 39001.25
$ws.Range("M116").Value = -35559.25
$ws.Range("H122").Value = 1435.8948
$ws.Range("I122").Value = 1413.25
$ws.Range("K122").Value = 4239.75
$ws.Range("M122").Value = -1789.75
$ws.Range("H128").Value = 141895
$ws.Range("J128").Value = 141895
$ws.Range("L128").Value = 141895
$ws.Range("N128").Value = -151855
$ws.Range("H130").Value = 116993
$ws.Range("J130").Value = 116993
$ws.Range("L130").Value = 116993
$ws.Range("N130").Value = -127033
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3187.86
$ws.Range("I122").Value = 2604.0857
$ws.Range("K122").Value = 7812.257100000001
$ws.Range("M122").Value = -5362.257100000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2717.6597
$ws.Range("I86").Value = 2683.7812
$ws.Range("J86").Value = 2789.9333
$ws.Range("K86").Value = 2683.7812
$ws.Range("L86").Value = 2789.9333
$ws.Range("M86").Value = -1560.7812
$ws.Range("N86").Value = -5035.933300000001
$ws.Range("H89").Value = 2717.6597
$ws.Range("I89").Value = 2683.7812
$ws.Range("J89").Value = 2789.9333
$ws.Range("K89").Value = 13418.906
$ws.Range("L89").Value = 13949.6665
$ws.Range("M89").Value = -7802.905999999999
$ws.Range("N89").Value = -25181.6665
$ws.Range("H107").Value = 2783.6316
$ws.Range("I107").Value = 2930.9375
$ws.Range("K107").Value = 2930.9375
$ws.Range("M107").Value = -1010.9375
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6139.4546
$ws.Range("J31").Value = 10798.6
$ws.Range("L31").Value = 10798.6
$ws.Range("N31").Value = -11388.6
$ws.Range("H34").Value = 6139.4546
$ws.Range("J34").Value = 10798.6
$ws.Range("L34").Value = 10798.6
$ws.Range("N34").Value = -11202.6
$ws.Range("H58").Value = 3346.0344
$ws.Range("I58").Value = 2985.0334
$ws.Range("K58").Value = 2985.0334
$ws.Range("M58").Value = -2782.0334
$ws.Range("H68").Value = 49966.555
$ws.Range("J68").Value = 49966.555
$ws.Range("L68").Value = 49966.555
$ws.Range("N68").Value = -51464.555
$ws.Range("H70").Value = 37450
$ws.Range("J70").Value = 37450
$ws.Range("L70").Value = 37450
$ws.Range("N70").Value = -38080
$ws.Range("H71").Value = 49966.555
$ws.Range("J71").Value = 49966.555
$ws.Range("L71").Value = 149899.665
$ws.Range("N71").Value = -157387.665
$ws.Range("H73").Value = 37450
$ws.Range("J73").Value = 37450
$ws.Range("L73").Value = 37450
$ws.Range("N73").Value = -39634
$ws.Range("H81").Value = 41328
$ws.Range("J81").Value = 41328
$ws.Range("L81").Value = 41328
$ws.Range("N81").Value = -43324
$ws.Range("H84").Value = 41328
$ws.Range("J84").Value = 41328
$ws.Range("L84").Value = 123984
$ws.Range("N84").Value = -133968
$ws.Range("H87").Value = 67960.8
$ws.Range("J87").Value = 67960.8
$ws.Range("L87").Value = 67960.8
$ws.Range("N87").Value = -70332.8
$ws.Range("H90").Value = 67960.8
$ws.Range("J90").Value = 67960.8
$ws.Range("L90").Value = 203882.4
$ws.Range("N90").Value = -215738.4
$ws.Range("H94").Value = 1739.4375
$ws.Range("J94").Value = 1526
$ws.Range("L94").Value = 1526
$ws.Range("N94").Value = -2428
$ws.Range("H134").Value = 1178.8334
$ws.Range("I134").Value = 1178.8334
$ws.Range("K134").Value = 3536.5002
$ws.Range("M134").Value = -1001.5002
$ws.Range("H136").Value = 3346.0344
$ws.Range("I136").Value = 2985.0334
$ws.Range("K136").Value = 8955.100199999999
$ws.Range("M136").Value = -6405.100199999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1490.4286
$ws.Range("I131").Value = 924.5714
$ws.Range("K131").Value = 2773.7142
$ws.Range("M131").Value = 2266.2858
$ws.Range("H132").Value = 3334665.8
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3334665.8
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = ""
$ws.Range("M132").Value = 30011992.2
$ws.Range("N132").Value = -30017052.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 390.6316
$ws.Range("I2").Value = 339.85715
$ws.Range("K2").Value = 339.85715
$ws.Range("M2").Value = -226.85715
$ws.Range("H13").Value = 999.5
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = ""
$ws.Range("H23").Value = 7666.6665
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 9000
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 9000
$ws.Range("M23").Value = -4777
$ws.Range("N23").Value = -9446
$ws.Range("H41").Value = 7057.25
$ws.Range("I41").Value = 6076.6665
$ws.Range("K41").Value = 6076.6665
$ws.Range("M41").Value = -5721.6665
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = ""
$ws.Range("N59").Value = 0
$ws.Range("H113").Value = 34349.65
$ws.Range("I113").Value = 6753.6665
$ws.Range("J113").Value = 75743.625
$ws.Range("K113").Value = 6753.6665
$ws.Range("L113").Value = 75743.625
$ws.Range("M113").Value = -4583.6665
$ws.Range("N113").Value = -80083.625
$ws.Range("H126").Value = 2034.6471
$ws.Range("I126").Value = 1509.3334
$ws.Range("J126").Value = 2625.625
$ws.Range("K126").Value = 4528.0002
$ws.Range("L126").Value = 7876.875
$ws.Range("M126").Value = -2058.0002
$ws.Range("N126").Value = -12816.875
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = ""
$ws.Range("N138").Value = 0
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3937.8635
$ws.Range("J61").Value = 7131.9
$ws.Range("L61").Value = 7131.9
$ws.Range("N61").Value = -7535.9
$ws.Range("H113").Value = 3937.8635
$ws.Range("J113").Value = 7131.9
$ws.Range("L113").Value = 7131.9
$ws.Range("N113").Value = -11471.9
